$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.200.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.072.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.39%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.070.81"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "

# Row 9
$ws.Range("E9").Value = "  +0.80%  "

# Row 10
$ws.Range("E10").Value = "  +1.97%  "

# Row 11
$ws.Range("E11").Value = "  -1.98%  "

# Row 12
$ws.Range("E12").Value = "  +3.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.576.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.257.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.85%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.074.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "

# Row 18
$ws.Range("E18").Value = "  +1.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.674"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("E27").Value = "  +2.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "

# Row 29
$ws.Range("E29").Value = "  +4.34%  "

# Row 30
$ws.Range("E30").Value = "  +0.28%  "

# Row 31
$ws.Range("E31").Value = "  +0.73%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "

# Row 33
$ws.Range("E33").Value = "  +3.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.78%  "

# Row 36
$ws.Range("E36").Value = "  +3.84%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "456.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.07%  "

# Row 39
$ws.Range("E39").Value = "  +2.83%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0822"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.953.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

# Row 43
$ws.Range("E43").Value = "  -1.93%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.35%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.21%  "

# Row 48
$ws.Range("E48").Value = "  +2.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0514"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.52%  "

# Row 51
$ws.Range("E51").Value = "  +0.43%  "
